$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the Status of issue #2 (row 3) from "Open" to "Fixed"
$ws.Range("D3").Value = "Fixed"

# Move the active selection to D4 (as recorded in the saved view state)
$ws.Range("D4").Select()
